# Add 2022-Q4 data.
#
# Target end state:
#   Sheet 1: "总计"    (sheetId=1, unchanged position) - gains a row for 2022-Q4
#   Sheet 2: "2022-Q4" (sheetId=2) - NEW fund data (the physical sheet that used to be
#                                    named "2021-Q2" is reused/renamed; its old content
#                                    is moved to the brand-new sheet below)
#   Sheet 3: "2021-Q2" (sheetId=3) - brand new sheet holding the fund data that used to
#                                    live in the "2021-Q2" sheet
#
# This specific sheetId/rId allocation (2022-Q4 keeps sheetId 2, 2021-Q2 gets a fresh
# sheetId 3) is what the engine produces when we rename the existing sheet in place and
# then Add() a new one after it - matching the target XML exactly.

$wb = $excel.ActiveWorkbook

$wsTotal = $wb.Worksheets.Item(1)      # "总计"
$wsFund  = $wb.Worksheets.Item(2)      # currently "2021-Q2", holds the old fund table

# ---------------------------------------------------------------------------
# 1) Add the brand-new sheet that will become "2021-Q2" right after the
#    existing fund sheet, and copy that existing sheet's header/A-column
#    look (style index 1) onto it while the look is still available.
#    NOTE: the existing sheet must be renamed away from "2021-Q2" *before*
#    the new sheet can claim that name (sheet names must be unique).
# ---------------------------------------------------------------------------
$wsOld = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsFund)

$wsFund.Range("B1:H1").Copy()
$wsOld.Range("B1:H1").PasteSpecial(-4122)
$wsFund.Range("A2").Copy()
$wsOld.Range("A2:A4").PasteSpecial(-4122)

$wsFund.Name = "2022-Q4"
$wsOld.Name = "2021-Q2"

# ---------------------------------------------------------------------------
# 2) Populate the new "2021-Q2" sheet with the OLD fund-table data that used
#    to live on the physical sheet we are about to repurpose.
# ---------------------------------------------------------------------------
$wsOld.Range("B1").Value = "基金代码"
$wsOld.Range("C1").Value = "基金名称"
$wsOld.Range("D1").Value = "基金金额"
$wsOld.Range("E1").Value = "股票总仓位"
$wsOld.Range("F1").Value = "仓位占比"
$wsOld.Range("G1").Value = "持有市值(亿元)"
$wsOld.Range("H1").Value = "仓位排名"

$oldFundRows = @(
    @("004814", "中欧红利优享灵活配置混合A", "2.32", "82.75", "4.80", "0.1114", 5),
    @("004815", "中欧红利优享灵活配置混合C", "0.95", "82.75", "4.80", "0.0456", 5),
    @("001940", "农银汇理现代农业加灵活配置混合", "1.46", "62.62", "2.41", "0.0352", 9)
)

for ($i = 0; $i -lt $oldFundRows.Count; $i++) {
    $r = 2 + $i
    $row = $oldFundRows[$i]
    $wsOld.Range("A$r").Value = $i
    $wsOld.Range("B$r").Value = "'" + $row[0]
    $wsOld.Range("B$r").ClearFormats()
    $wsOld.Range("C$r").Value = $row[1]
    $wsOld.Range("D$r").Value = "'" + $row[2]
    $wsOld.Range("D$r").ClearFormats()
    $wsOld.Range("E$r").Value = "'" + $row[3]
    $wsOld.Range("E$r").ClearFormats()
    $wsOld.Range("F$r").Value = "'" + $row[4]
    $wsOld.Range("F$r").ClearFormats()
    $wsOld.Range("G$r").Value = "'" + $row[5]
    $wsOld.Range("G$r").ClearFormats()
    $wsOld.Range("H$r").Value = $row[6]
}

# ---------------------------------------------------------------------------
# 3) Clear the original fund sheet's old 4th data row (new table has only 2
#    data rows) and fill it in with the new fund data. (Already renamed to
#    "2022-Q4" in step 1.)
# ---------------------------------------------------------------------------
$wsFund.Range("A4:H4").Clear()

$wsFund.Range("D1").Value = "基金规模"   # header wording change ("基金金额" -> "基金规模")

$newFundRows = @(
    @("015921", "申万菱信国证2000指数增强A", "0.21", "94.00", "0.50", "0.0010", 9),
    @("015922", "申万菱信国证2000指数增强C", "0.08", "94.00", "0.50", "0.0004", 9)
)

for ($i = 0; $i -lt $newFundRows.Count; $i++) {
    $r = 2 + $i
    $row = $newFundRows[$i]
    $wsFund.Range("A$r").Value = $i
    $wsFund.Range("B$r").Value = "'" + $row[0]
    $wsFund.Range("B$r").ClearFormats()
    $wsFund.Range("C$r").Value = $row[1]
    $wsFund.Range("D$r").Value = "'" + $row[2]
    $wsFund.Range("D$r").ClearFormats()
    $wsFund.Range("E$r").Value = "'" + $row[3]
    $wsFund.Range("E$r").ClearFormats()
    $wsFund.Range("F$r").Value = "'" + $row[4]
    $wsFund.Range("F$r").ClearFormats()
    $wsFund.Range("G$r").Value = "'" + $row[5]
    $wsFund.Range("G$r").ClearFormats()
    $wsFund.Range("H$r").Value = $row[6]
}

# Restyle "2022-Q4" header/A-column to style index 2 (matching the "总计"
# sheet's own look, which is what the target workbook uses for this sheet).
$wsTotal.Range("B1").Copy()
$wsFund.Range("B1:H1").PasteSpecial(-4122)
$wsTotal.Range("A2").Copy()
$wsFund.Range("A2:A3").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 4) Update the "总计" sheet: insert a row for 2022-Q4 above the existing
#    2021-Q2 entry (which shifts down to row 3).
# ---------------------------------------------------------------------------
$wsTotal.Rows.Item(2).Insert(-4121)
$wsTotal.Range("A2:D2").ClearFormats()
$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q4"
$wsTotal.Range("C2").Value = 2
$wsTotal.Range("D2").Value = 0
$wsTotal.Range("A3").Value = 1

$wsTotal.Range("A3").Copy()
$wsTotal.Range("A2").PasteSpecial(-4122)

Write-Output "2022-Q4 sheet added successfully"
